$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.629.11'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').Value = '2.587.74'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '508.35'
$ws.Range('E5').Value = '  -1.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.08'
$ws.Range('E6').Value = '  -3.95%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('E8').Value = '  -6.00%  '
$ws.Range('D9').Value = '2.594.30'
$ws.Range('E9').Value = '  -3.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.62'
$ws.Range('E10').Value = '  +6.13%  '
$ws.Range('E11').Value = '  -3.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.348'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').Value = '3.039.25'
$ws.Range('E14').Value = '  -1.70%  '
$ws.Range('D15').Value = '60.604.36'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.75'
$ws.Range('E16').Value = '  -4.23%  '
$ws.Range('E17').Value = '  -1.85%  '
$ws.Range('D18').Value = '2.593.51'
$ws.Range('E18').Value = '  -2.89%  '
$ws.Range('E19').Value = '  -1.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '347.51'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.53'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.14'
$ws.Range('E22').Value = '  -2.25%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.43'
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('E26').Value = '  -2.12%  '
$ws.Range('D27').Value = '2.697.34'
$ws.Range('E27').Value = '  -3.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  +0.49%  '
$ws.Range('D29').Value = '0.0₃0850'
$ws.Range('E29').Value = '  -3.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.43'
$ws.Range('E30').Value = '  -3.36%  '
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.42'
$ws.Range('E32').Value = '  -2.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '152.92'
$ws.Range('E33').Value = '  -3.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.57'
$ws.Range('E34').Value = '  -2.50%  '
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('E37').Value = '  -3.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.855'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.49'
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.850'
$ws.Range('E40').Value = '  -4.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.28'
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '298.41'
$ws.Range('E43').Value = '  -2.03%  '
$ws.Range('E44').Value = '  -3.61%  '
$ws.Range('E45').Value = '  -2.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0561'
$ws.Range('E46').Value = '  -4.70%  '
$ws.Range('E47').Value = '  +0.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.81'
$ws.Range('E48').Value = '  -3.01%  '
$ws.Range('E49').Value = '  -3.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0235'
$ws.Range('E50').Value = '  -2.20%  '
$ws.Range('E51').Value = '  +0.20%  '

Write-Host "Updated cryptos list"